$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "67.484.76", "0.999").
# Force Text format first so Excel does not auto-convert these numeric-looking
# strings into actual numbers when we assign them below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.484.76"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "3.321.19"
$ws.Range("E3").Value = "  +1.95%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "575.99"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").Value = "174.08"
$ws.Range("E6").Value = "  -2.18%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("E8").Value = "  +1.48%  "

$ws.Range("D9").Value = "3.318.98"
$ws.Range("E9").Value = "  +2.08%  "

$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  +1.98%  "

$ws.Range("D11").Value = "0.579"
$ws.Range("E11").Value = "  +1.52%  "

$ws.Range("D12").Value = "45.83"
$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("D14").Value = "704.19"
$ws.Range("E14").Value = "  +3.78%  "

$ws.Range("D15").Value = "3.858.10"
$ws.Range("E15").Value = "  +1.96%  "

$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  +1.29%  "

$ws.Range("D17").Value = "67.488.80"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "3.325.09"
$ws.Range("E19").Value = "  +1.71%  "

$ws.Range("D20").Value = "17.35"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "10.96"
$ws.Range("E21").Value = "  +2.57%  "

$ws.Range("D22").Value = "0.890"
$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("D23").Value = "5.37"
$ws.Range("E23").Value = "  +5.00%  "

$ws.Range("D24").Value = "16.86"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").Value = "98.66"
$ws.Range("E25").Value = "  +0.61%  "

$ws.Range("D26").Value = "3.85"
$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("D27").Value = "2.67"
$ws.Range("E27").Value = "  -1.77%  "

$ws.Range("D28").Value = "9.31"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "33.25"
$ws.Range("E29").Value = "  +2.35%  "

$ws.Range("D30").Value = "8.49"
$ws.Range("E30").Value = "  +1.34%  "

$ws.Range("E31").Value = "  +7.20%  "

$ws.Range("D32").Value = "568.15"
$ws.Range("E32").Value = "  -2.10%  "

$ws.Range("D33").Value = "10.95"
$ws.Range("E33").Value = "  +1.61%  "

$ws.Range("D34").Value = "0.105"
$ws.Range("E34").Value = "  +1.63%  "

$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("D36").Value = "3.692.35"
$ws.Range("E36").Value = "  -4.31%  "

$ws.Range("D37").Value = "56.88"
$ws.Range("E37").Value = "  +3.20%  "

$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  -1.67%  "

$ws.Range("D39").Value = "34.16"
$ws.Range("E39").Value = "  +6.85%  "

$ws.Range("E40").Value = "  +0.89%  "

$ws.Range("D41").Value = "3.15"
$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("D42").Value = "2.61"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -4.46%  "

$ws.Range("D46").Value = "0.0403"
$ws.Range("E46").Value = "  -1.30%  "

$ws.Range("D47").Value = "2.68"
$ws.Range("E47").Value = "  +7.83%  "

$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("E50").Value = "  -4.77%  "

$ws.Range("D51").Value = "128.95"
$ws.Range("E51").Value = "  -0.54%  "

# Row 43 and row 44 swap coins: PEPE now ranks above TheGraph, with refreshed
# price/volume figures for both.
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0667"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.333"
$ws.Range("E44").Value = "  +1.37%  "
